# Pushing the new code
#
# The "Couponid" sheet (sheet2) becomes the active/selected tab, its A2 cell
# is changed from the number 820009 to the text "7676735387a" (added as a
# new shared string), and its selection moves to A2. Consequently the
# previously-active "CardDetails" sheet (sheet3) is no longer the selected
# tab (its own selection stays at G4, untouched).

$wb = $excel.ActiveWorkbook

$couponSheet = $wb.Worksheets.Item("Couponid")

# Make Couponid the active sheet (updates workbook-level activeTab and
# moves tabSelected between the worksheets' sheetViews).
$couponSheet.Activate()

# Replace the numeric value with the new text value.
$couponSheet.Range("A2").Value = "7676735387a"

# Update the on-sheet selection to A2.
$couponSheet.Range("A2").Select()
